$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Investor" column (E) with distribution/investor info per commitment %
$ws.Range("E1").Value = "Investor"
$ws.Range("E2").Value = "Kalaari Capital"
$ws.Range("E3").Value = "Accel"

# Update active cell selection to E4 as left by the editor
$ws.Range("E4").Select()
